# Added transaction and idempotency to PortfolioAttribution
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update investment dates for rows 3, 4, and 6
$ws.Range("C3").Value = Get-Date -Year 2021 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("C4").Value = Get-Date -Year 2022 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("C6").Value = Get-Date -Year 2023 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0

# Update quantity for row 6
$ws.Range("E6").Value = -3000

# Column F width (best-fit) - set based on autofit after content width needs
$ws.Columns("F").ColumnWidth = 14.9296875

# Update the active selection to L24
$ws.Range("L24").Select()
